$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header) values tweaked
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON): B2, D2, E2 deleted entirely; C2 value replaced
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -0.7960953803542401
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 (STR): values replaced
$ws.Range("B3").Value = -0.8956985217115051
$ws.Range("C3").Value = -0.22812556512210955
$ws.Range("D3").Value = -1.6748434028007984
$ws.Range("E3").Value = 2.1048551030717273

# Selection narrowed to the edited block
$ws.Range("B1:E3").Select() | Out-Null
